$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2548.4
$ws.Range("C2").Value = 2495
$ws.Range("D2").Value = 2523
$ws.Range("E2").Value = 2531.1999999999998
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 2514.0500000000002
$ws.Range("B3").Value = 450.75
$ws.Range("C3").Value = 442.75
$ws.Range("D3").Value = 450
$ws.Range("E3").Value = 449.95
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 444.85
$ws.Range("B4").Value = 1711.7
$ws.Range("C4").Value = 1685.9
$ws.Range("D4").Value = 1694
$ws.Range("E4").Value = 1696.25
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 1701
$ws.Range("B5").Value = 7420
$ws.Range("C5").Value = 7315.05
$ws.Range("D5").Value = 7405
$ws.Range("E5").Value = 7407.65
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 7357.65
$ws.Range("B6").Value = 234.85
$ws.Range("C6").Value = 229.25
$ws.Range("D6").Value = 234.5
$ws.Range("E6").Value = 234.3
$ws.Range("F6").Value = 141
$ws.Range("G6").Value = 230.85
$ws.Range("B7").Value = 210
$ws.Range("C7").Value = 205.75
$ws.Range("D7").Value = 209
$ws.Range("E7").Value = 209.05
$ws.Range("F7").Value = 359
$ws.Range("G7").Value = 206.3
$ws.Range("B8").Value = 360
$ws.Range("C8").Value = 351.05
$ws.Range("D8").Value = 354.5
$ws.Range("E8").Value = 354.7
$ws.Range("F8").Value = 227
$ws.Range("G8").Value = 352.1
$ws.Range("B9").Value = 652
$ws.Range("C9").Value = 632.79999999999995
$ws.Range("D9").Value = 649.15
$ws.Range("E9").Value = 649.25
$ws.Range("F9").Value = 45
$ws.Range("G9").Value = 636.45000000000005
$ws.Range("B10").Value = 4200
$ws.Range("C10").Value = 3931.3
$ws.Range("D10").Value = 4181
$ws.Range("E10").Value = 4180.3500000000004
$ws.Range("F10").Value = 28
$ws.Range("G10").Value = 3940
$ws.Range("B11").Value = 155.19999999999999
$ws.Range("C11").Value = 151.05000000000001
$ws.Range("D11").Value = 154.1
$ws.Range("E11").Value = 154.05000000000001
$ws.Range("F11").Value = 217
$ws.Range("G11").Value = 151.55000000000001
$ws.Range("B12").Value = 1344.35
$ws.Range("C12").Value = 1331
$ws.Range("D12").Value = 1335
$ws.Range("E12").Value = 1337.95
$ws.Range("F12").Value = 15
$ws.Range("G12").Value = 1339.55
$ws.Range("B13").Value = 1612
$ws.Range("C13").Value = 1581.3
$ws.Range("D13").Value = 1609.15
$ws.Range("E13").Value = 1609.4
$ws.Range("F13").Value = 253
$ws.Range("G13").Value = 1581.85
$ws.Range("B14").Value = 526.65
$ws.Range("C14").Value = 516.54999999999995
$ws.Range("D14").Value = 519.25
$ws.Range("E14").Value = 519.25
$ws.Range("G14").Value = 521
$ws.Range("B15").Value = 994.85
$ws.Range("C15").Value = 966.85
$ws.Range("D15").Value = 991
$ws.Range("E15").Value = 991.2
$ws.Range("F15").Value = 331
$ws.Range("G15").Value = 972.6
$ws.Range("B16").Value = 1517.55
$ws.Range("C16").Value = 1479.65
$ws.Range("D16").Value = 1515.55
$ws.Range("E16").Value = 1514.65
$ws.Range("F16").Value = 42
$ws.Range("G16").Value = 1483.15
$ws.Range("B17").Value = 1471.9
$ws.Range("C17").Value = 1449.3
$ws.Range("D17").Value = 1464.9
$ws.Range("E17").Value = 1464.35
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 1467.25
$ws.Range("B18").Value = 694.05
$ws.Range("C18").Value = 678.45
$ws.Range("D18").Value = 680.85
$ws.Range("E18").Value = 681.5
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 693.3
$ws.Range("B19").Value = 510.6
$ws.Range("C19").Value = 497.1
$ws.Range("D19").Value = 504.95
$ws.Range("E19").Value = 507.55
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 498.95
$ws.Range("B20").Value = 1675
$ws.Range("C20").Value = 1639.45
$ws.Range("D20").Value = 1655
$ws.Range("E20").Value = 1656.85
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 1659.45
$ws.Range("B21").Value = 281.35000000000002
$ws.Range("C21").Value = 276.10000000000002
$ws.Range("D21").Value = 278.05
$ws.Range("E21").Value = 277.7
$ws.Range("F21").Value = 27
$ws.Range("G21").Value = 277.7
$ws.Range("B22").Value = 277.35000000000002
$ws.Range("C22").Value = 273.25
$ws.Range("D22").Value = 274.89999999999998
$ws.Range("E22").Value = 274.8
$ws.Range("F22").Value = 204
$ws.Range("G22").Value = 274.35000000000002
$ws.Range("B23").Value = 2426.1999999999998
$ws.Range("C23").Value = 2398.6
$ws.Range("D23").Value = 2421
$ws.Range("E23").Value = 2420.1999999999998
$ws.Range("F23").Value = 77
$ws.Range("G23").Value = 2418.6
$ws.Range("B24").Value = 596
$ws.Range("C24").Value = 584.54999999999995
$ws.Range("D24").Value = 595.9
$ws.Range("E24").Value = 594.70000000000005
$ws.Range("F24").Value = 260
$ws.Range("G24").Value = 585.04999999999995
$ws.Range("B25").Value = 680
$ws.Range("C25").Value = 673.8
$ws.Range("D25").Value = 674
$ws.Range("E25").Value = 676.65
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 675
$ws.Range("B26").Value = 978.8
$ws.Range("C26").Value = 971
$ws.Range("D26").Value = 971
$ws.Range("E26").Value = 973.05
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 973.65
$ws.Range("B27").Value = 709.9
$ws.Range("C27").Value = 704.8
$ws.Range("D27").Value = 706.1
$ws.Range("E27").Value = 705.6
$ws.Range("F27").Value = 93
$ws.Range("G27").Value = 707.35
$ws.Range("B28").Value = 283.89999999999998
$ws.Range("C28").Value = 278.89999999999998
$ws.Range("D28").Value = 281.89999999999998
$ws.Range("E28").Value = 281.55
$ws.Range("F28").Value = 288
$ws.Range("G28").Value = 281.14999999999998
$ws.Range("B29").Value = 131.94999999999999
$ws.Range("C29").Value = 130.4
$ws.Range("D29").Value = 131
$ws.Range("E29").Value = 131
$ws.Range("F29").Value = 412
$ws.Range("G29").Value = 131.35
$ws.Range("B30").Value = 9340.7999999999993
$ws.Range("C30").Value = 9090
$ws.Range("D30").Value = 9295
$ws.Range("E30").Value = 9317.75
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 9317.75

$ws.Range("G12").Select()
